# "added invalid data driven script"
# Add a new worksheet "InvalidLogin" right after "ValidLogin" that holds a
# negative/invalid login test-data row, make it the active sheet/tab, and
# bump its zoom level.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after ValidLogin (Before=$null, After=$ws1).
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "InvalidLogin"

# Header row, matching the existing ValidLogin layout.
$new.Range("A1").Value = "UserName"
$new.Range("B1").Value = "Password"

# Invalid credentials data row.
$new.Range("A2").Value = "abc"
$new.Range("B2").Value = "xyz"

# Match the selection/view state captured in the target sheet.
$new.Range("A3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220
